# PROJETO_COMISSOES_V2 - separar calculo de comissoes por recebimento
# Ajustes na planilha "Plan1":
#  - S2 deixa de ser o texto "048341" e passa a ser o numero 48341
#  - A4 deixa de ser o texto "CS-Test" e passa a ser o numero 999997
#  - Larguras de colunas (F:K e Q) recalculadas (auto-ajuste)
#  - Selecao ativa passa para F8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numero NF (S2): era texto "048341", agora numero 48341 ---
$ws.Range("S2").Value = 48341

# --- Processo (A4): era texto "CS-Test", agora numero 999997 ---
$ws.Range("A4").ClearFormats()
$ws.Range("A4").Value = 999997

# --- Reajuste das larguras de coluna (auto fit) ---
# Os valores abaixo foram escolhidos para aproximar, o mais perto possivel, as
# larguras finais "bestFit" originais (14.140625 / 8.85546875 / 16.5703125 /
# 17.7109375 / 13.85546875 / 15.7109375 / 25.28515625), compensando o
# arredondamento interno do motor de largura de coluna.
$ws.Columns("F").ColumnWidth = 13.307291666666666
$ws.Columns("G").ColumnWidth = 8.022135416666666
$ws.Columns("H").ColumnWidth = 15.736979166666666
$ws.Columns("I").ColumnWidth = 16.877604166666668
$ws.Columns("J").ColumnWidth = 13.022135416666666
$ws.Columns("K").ColumnWidth = 14.877604166666666
$ws.Columns("Q").ColumnWidth = 24.451822916666668

# --- Selecao final do usuario ---
$ws.Range("F8").Select()
